$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = "[-, -, Leandro-Mec. Manut. Equip. Ind., -]"
$ws.Range("D2").Value = "[-, Leandro-M. S. Ar Cond., -, -]"
$ws.Range("E2").Value = "[Aderci-Fresagem, Ivan-Trat. Térmicos, Victor-Usin. CNC, Rogério-Retificação]"
$ws.Range("F2").Value = "[Victor-Usin. CNC, Aderci-Fresagem, Rogério-Retificação, Ivan-Trat. Térmicos]"

# Row 3
$ws.Range("B3").Value = "-"
$ws.Range("D3").Value = "[-, Leandro-M. S. Ar Cond., -, -]"
$ws.Range("E3").Value = "[Aderci-Fresagem, Valmir-Calderaria, Victor-Usin. CNC, Rogério-Retificação]"
$ws.Range("F3").Value = "[Ludoff-Comam. Pneumáticos, Aderci-Fresagem, Rogério-Retificação, Ivan-Trat. Térmicos]"

# Row 4
$ws.Range("B4").Value = "-"
$ws.Range("C4").Value = "Cleison-Elem"
$ws.Range("D4").Value = "[Leandro-Mec. Manut. Equip. Ind., -, -, -]"
$ws.Range("E4").Value = "[Gisele-Ens. Dest. não Dest., Valmir-Calderaria, Anderson J.-M. A. Comp; Cad / CAM, Humberto-Cont. Lóg. Prog. CLP]"
$ws.Range("F4").Value = "[Ludoff-Coman. Hidráulicos, Anselmo-M. Motor Endot., Leandro-M. S. Ar Cond., Ivan-Trat. Térmicos]"

# Row 6
$ws.Range("B6").Value = "-"
$ws.Range("C6").Value = "Cleison-Elem"
$ws.Range("D6").Value = "[Leandro-Mec. Manut. Equip. Ind., -, Vinicius-Metrologia 2, Ludoff-Coman. Hidráulicos]"
$ws.Range("E6").Value = "[Gisele-Ens. Dest. não Dest., Valmir-Calderaria, Anderson J.-M. A. Comp; Cad / CAM, Humberto-Cont. Lóg. Prog. CLP]"
$ws.Range("F6").Value = "[Ivan-Tec. Soldagem, Anselmo-M. Motor Endot., Ludoff-Comam. Pneumáticos, Joel Lima-Tec. Fundição]"

# Row 7
$ws.Range("B7").Value = "-"
$ws.Range("C7").Value = "[-, Vinicius-Metrologia 2, -, -]"
$ws.Range("D7").Value = "[-, Leandro-M. S. Ar Cond., Vinicius-Metrologia 2, Ludoff-Coman. Hidráulicos]"
$ws.Range("E7").Value = "[Gisele-Ens. Dest. não Dest., Valmir-Calderaria, Anderson J.-M. A. Comp; Cad / CAM, Humberto-Cont. Lóg. Prog. CLP]"
$ws.Range("F7").Value = "[Ivan-Tec. Soldagem, Anselmo-M. Motor Endot., Ludoff-Comam. Pneumáticos, Joel Lima-Tec. Fundição]"

# Row 8
$ws.Range("B8").Value = "-"
$ws.Range("C8").Value = "[-, Vinicius-Metrologia 2, -, -]"
$ws.Range("D8").Value = "[Ivan-Tec. Soldagem, Joel Lima-Tec. Fundição, Leandro-Mec. Manut. Equip. Ind., Ludoff-Coman. Hidráulicos]"
$ws.Range("E8").Value = "[Gisele-Ens. Dest. não Dest., Victor-Usin. CNC, Anderson J.-M. A. Comp; Cad / CAM, Humberto-Cont. Lóg. Prog. CLP]"
$ws.Range("F8").Value = "[Ivan-Tec. Soldagem, Anselmo-M. Motor Endot., Ludoff-Comam. Pneumáticos, Joel Lima-Tec. Fundição]"
